# adicionando alteracoes perdidas pelo merge
#
# This script reconstructs the missing merge changes:
#  1. Rows 81-84: column F ("Data") was stored as text ("15/09/2025"); it
#     should actually be a real date serial (45915) formatted the same way
#     as the rest of column F (YYYY-MM-DD HH:MM:SS).
#  2. Four new "Farmacia" rows with time 11:30:08 and a numeric date (85-88).
#  3. Four more new "Farmacia" rows with time 11:31:33 and the date kept as
#     text "15/09/2025" (89-92), matching how rows 81-84 originally looked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number format used by the existing date column (style index already in
# use for F2:F80) - reuse it so no new style gets created.
$dateFormat = $ws.Range("F2").NumberFormat

# --- 1. Fix rows 81-84: turn the textual date into a real numeric date ---
for ($r = 81; $r -le 84; $r++) {
    $ws.Cells.Item($r, 6).Value = 45915
    $ws.Cells.Item($r, 6).NumberFormat = $dateFormat
}

# --- 2 & 3. Append the new "Farmacia" rows ---
$funcionario = "admin"
$insumos = @("seringa", "algodão", "gazes", "luvas")
$consumos = @(23, 333, 44, 55)
$setor = "Farmácia"

# Rows 85-88: hora 11:30:08, data numérica 45915
$row = 85
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($row, 1).Value = $funcionario
    $ws.Cells.Item($row, 2).Value = $insumos[$i]
    $ws.Cells.Item($row, 3).Value = $consumos[$i]
    $ws.Cells.Item($row, 4).Value = $setor
    $ws.Cells.Item($row, 5).Value = "11:30:08"
    $ws.Cells.Item($row, 6).Value = 45915
    $ws.Cells.Item($row, 6).NumberFormat = $dateFormat
    $row++
}

# Rows 89-92: hora 11:31:33, data textual "15/09/2025"
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($row, 1).Value = $funcionario
    $ws.Cells.Item($row, 2).Value = $insumos[$i]
    $ws.Cells.Item($row, 3).Value = $consumos[$i]
    $ws.Cells.Item($row, 4).Value = $setor
    $ws.Cells.Item($row, 5).Value = "11:31:33"
    $ws.Cells.Item($row, 6).Value = "15/09/2025"
    $row++
}
